$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "24/10/2025"
$ws.Range("B14").Value = "Motor Lublin"
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "Widzew Lodz"
$ws.Range("F14").Value = "L"
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 0.51
$ws.Range("L14").Value = 2.65
$ws.Range("M14").Value = 13
$ws.Range("N14").Value = 16
$ws.Range("O14").Value = 4
$ws.Range("P14").Value = 5
